$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of summary-statistic data to append below the existing table (rows 5-14)
$data = @(
    @("wm", "m", 321, "median depth change", "sum of change of median depth per week is << 0"),
    @("sm", "m", 321, "median depth change", "sum of absolute change of median depth per week >>0"),
    @("sr", "m", 321, "median depth change", "sum of change of median depth per week == 0"),
    @("sr", "m", 321, "median depth change", "mean change of median depth per week == 0"),
    @("wr", "m", 321, "median depth change", "sd of change of median depth per week >> 0"),
    @("sr", "m", 321, "median depth change", "sd of change of median depth per week == 0"),
    @("sr", "f", 308, "median depth change", "sum of change of median depth per week == 0"),
    @("sr", "f", 308, "median depth change", "mean change of median depth per week == 0"),
    @("wr", "f", 308, "vertical movement", "predominant vertical movement per week = DVM"),
    @("wm", "f", 308, "median depth change", "sum of change of median depth per week is << 0")
)

$startRow = 5

# Column A is filled in first for every new row (as the original author did),
# which determines the order new strings were interned into the shared
# strings table; then the remaining columns B-E are filled in row by row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$ws.Range("D11").Select()
